# Semana 6, Clases 11 y 12
# Rename "Ejemplo Bomba" sheet to "Ejemplo Sistema Bombeo" (Solver-linked
# defined names that point at this sheet are updated automatically by the
# rename).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ejemplo Bomba")
$ws.Name = "Ejemplo Sistema Bombeo"

# Re-enter the "Valor <= Limite" check formula across I15:I22 in one shot so
# it collapses into a single shared formula (same behaviour as typing the
# formula once in I15 and filling it down through I22).
$ws.Range("I15:I22").Formula = "=F15<=H15"
